$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing row (20) down into the new row (21)
# so the new row matches the established look (number formats, styles, etc.)
$ws.Range("A20:H20").Copy()
$ws.Range("A21:H21").PasteSpecial(-4122)  # xlPasteFormats

# New weigh-in entry added 14/05/2018
$ws.Range("A21").Value2 = 20
$ws.Range("B21").Value2 = 43234
$ws.Range("C21").Value2 = 14.1
$ws.Range("D21").Value2 = 93.4
$ws.Range("E21").Value2 = 205.9
$ws.Range("F21").Formula = "=E21-E20"
$ws.Range("G21").Formula = "=ROUND((D21/1.88)/1.88,2)"
$ws.Range("H21").Value2 = 20.7
